$wb = $excel.ActiveWorkbook

# --- Step1_Data sheet: raw probability values ---
$ws1 = $wb.Worksheets.Item("Step1_Data")
$ws1.Cells.Item(2, 4).Value = 0.05027127526523962
$ws1.Cells.Item(2, 5).Value = 0.1222581803075866
$ws1.Cells.Item(2, 6).Value = 0.4464464645675328
$ws1.Cells.Item(2, 7).Value = 0.1221089397079538
$ws1.Cells.Item(2, 8).Value = 0.01700147511216995
$ws1.Cells.Item(2, 9).Value = 0.01432126311808001
$ws1.Cells.Item(2, 10).Value = 0
$ws1.Cells.Item(2, 11).Value = 0.006529660145584581
$ws1.Cells.Item(2, 12).Value = 0.00108249413915575
$ws1.Cells.Item(2, 13).Value = 0
$ws1.Cells.Item(2, 14).Value = 0.0009129630852946645
$ws1.Cells.Item(2, 15).Value = 0.007392084423924092
$ws1.Cells.Item(2, 16).Value = 0.05450611437748384
$ws1.Cells.Item(2, 17).Value = 0
$ws1.Cells.Item(2, 18).Value = 0
$ws1.Cells.Item(2, 19).Value = 0
$ws1.Cells.Item(2, 20).Value = 0
$ws1.Cells.Item(2, 21).Value = 0
$ws1.Cells.Item(2, 22).Value = 0
$ws1.Cells.Item(2, 23).Value = 0.02601491857744207
$ws1.Cells.Item(2, 24).Value = 0
$ws1.Cells.Item(2, 25).Value = 0.02007506437053477
$ws1.Cells.Item(2, 26).Value = 0
$ws1.Cells.Item(2, 27).Value = 0.03870125088835333
$ws1.Cells.Item(2, 28).Value = 0.01350321931201975
$ws1.Cells.Item(2, 29).Value = 0.02031645418555763
$ws1.Cells.Item(2, 30).Value = 0
$ws1.Cells.Item(2, 31).Value = 0.02690872026750446
$ws1.Cells.Item(2, 32).Value = 0
$ws1.Cells.Item(2, 33).Value = 0.01164945814858224
$ws1.Cells.Item(2, 34).Value = 0
$ws1.Cells.Item(2, 35).Value = 0.0
$ws1.Cells.Item(3, 4).Value = 0.05026589338730262
$ws1.Cells.Item(3, 5).Value = 0.1590547699181399
$ws1.Cells.Item(3, 6).Value = 0.3563179814252341
$ws1.Cells.Item(3, 7).Value = 0.1745522467868859
$ws1.Cells.Item(3, 8).Value = 0.0224624715379248
$ws1.Cells.Item(3, 9).Value = 0.02200306747644844
$ws1.Cells.Item(3, 10).Value = 0
$ws1.Cells.Item(3, 11).Value = 0.007983512738412414
$ws1.Cells.Item(3, 12).Value = 0
$ws1.Cells.Item(3, 13).Value = 0
$ws1.Cells.Item(3, 14).Value = 0
$ws1.Cells.Item(3, 15).Value = 0
$ws1.Cells.Item(3, 16).Value = 0.02799846500541998
$ws1.Cells.Item(3, 17).Value = 0
$ws1.Cells.Item(3, 18).Value = 0
$ws1.Cells.Item(3, 19).Value = 0
$ws1.Cells.Item(3, 20).Value = 0
$ws1.Cells.Item(3, 21).Value = 0
$ws1.Cells.Item(3, 22).Value = 0
$ws1.Cells.Item(3, 23).Value = 0.02974344171050223
$ws1.Cells.Item(3, 24).Value = 0
$ws1.Cells.Item(3, 25).Value = 0.02506773408396289
$ws1.Cells.Item(3, 26).Value = 0
$ws1.Cells.Item(3, 27).Value = 0.03663759912565145
$ws1.Cells.Item(3, 28).Value = 0.01435563589175506
$ws1.Cells.Item(3, 29).Value = 0.02150125693530683
$ws1.Cells.Item(3, 30).Value = 0
$ws1.Cells.Item(3, 31).Value = 0.03487617479415477
$ws1.Cells.Item(3, 32).Value = 0
$ws1.Cells.Item(3, 33).Value = 0.01717974918289883
$ws1.Cells.Item(3, 34).Value = 0
$ws1.Cells.Item(3, 35).Value = 0.0
$ws1.Cells.Item(4, 4).Value = 0.2943213283311363
$ws1.Cells.Item(4, 5).Value = 0.05230181659181025
$ws1.Cells.Item(4, 6).Value = 0.3597687889456396
$ws1.Cells.Item(4, 7).Value = 0
$ws1.Cells.Item(4, 8).Value = 0.007183637722097465
$ws1.Cells.Item(4, 9).Value = 0
$ws1.Cells.Item(4, 10).Value = 0
$ws1.Cells.Item(4, 11).Value = 0
$ws1.Cells.Item(4, 12).Value = 0.02926172633794116
$ws1.Cells.Item(4, 13).Value = 0.0007519634961031186
$ws1.Cells.Item(4, 14).Value = 0.003963267529695806
$ws1.Cells.Item(4, 15).Value = 0.09286479354853734
$ws1.Cells.Item(4, 16).Value = 0.02807313113895594
$ws1.Cells.Item(4, 17).Value = 0
$ws1.Cells.Item(4, 18).Value = 0
$ws1.Cells.Item(4, 19).Value = 0
$ws1.Cells.Item(4, 20).Value = 0
$ws1.Cells.Item(4, 21).Value = 0
$ws1.Cells.Item(4, 22).Value = 0.05632289886578625
$ws1.Cells.Item(4, 23).Value = 0
$ws1.Cells.Item(4, 24).Value = 0
$ws1.Cells.Item(4, 25).Value = 0
$ws1.Cells.Item(4, 26).Value = 0
$ws1.Cells.Item(4, 27).Value = 0.03898213655671936
$ws1.Cells.Item(4, 28).Value = 0
$ws1.Cells.Item(4, 29).Value = 0.005305971531258637
$ws1.Cells.Item(4, 30).Value = 0.01304209320757299
$ws1.Cells.Item(4, 31).Value = 0
$ws1.Cells.Item(4, 32).Value = 0.0137459660589164
$ws1.Cells.Item(4, 33).Value = 0.004110480137829569
$ws1.Cells.Item(4, 34).Value = 0
$ws1.Cells.Item(4, 35).Value = 0.0
$ws1.Cells.Item(5, 4).Value = 0.1729002386830045
$ws1.Cells.Item(5, 5).Value = 0.0076487841857635
$ws1.Cells.Item(5, 6).Value = 0.5570038134192203
$ws1.Cells.Item(5, 7).Value = 0.04274474983070176
$ws1.Cells.Item(5, 8).Value = 0.0413738049007074
$ws1.Cells.Item(5, 9).Value = 0.004144928657331091
$ws1.Cells.Item(5, 10).Value = 0
$ws1.Cells.Item(5, 11).Value = 0.01854448849737875
$ws1.Cells.Item(5, 12).Value = 0
$ws1.Cells.Item(5, 13).Value = 0
$ws1.Cells.Item(5, 14).Value = 0
$ws1.Cells.Item(5, 15).Value = 0.02419996875679432
$ws1.Cells.Item(5, 16).Value = 0.02979261811386093
$ws1.Cells.Item(5, 17).Value = 0
$ws1.Cells.Item(5, 18).Value = 0
$ws1.Cells.Item(5, 19).Value = 0
$ws1.Cells.Item(5, 20).Value = 0
$ws1.Cells.Item(5, 21).Value = 0
$ws1.Cells.Item(5, 22).Value = 0.001474395082352065
$ws1.Cells.Item(5, 23).Value = 0.01600436991989613
$ws1.Cells.Item(5, 24).Value = 0
$ws1.Cells.Item(5, 25).Value = 0.01544338806221371
$ws1.Cells.Item(5, 26).Value = 0
$ws1.Cells.Item(5, 27).Value = 0.03615203805171234
$ws1.Cells.Item(5, 28).Value = 0.003994287981948967
$ws1.Cells.Item(5, 29).Value = 0.01867782648321598
$ws1.Cells.Item(5, 30).Value = 0
$ws1.Cells.Item(5, 31).Value = 0.009900299373898488
$ws1.Cells.Item(5, 32).Value = 0
$ws1.Cells.Item(5, 33).Value = 0.0
$ws1.Cells.Item(5, 34).Value = 0.0
$ws1.Cells.Item(5, 35).Value = 0.0
$ws1.Cells.Item(6, 4).Value = 0.03304963354092785
$ws1.Cells.Item(6, 5).Value = 0.1714022291684707
$ws1.Cells.Item(6, 6).Value = 0.3412580169285142
$ws1.Cells.Item(6, 7).Value = 0.1657837823912913
$ws1.Cells.Item(6, 8).Value = 0.0467640209035246
$ws1.Cells.Item(6, 9).Value = 0.03306616298755186
$ws1.Cells.Item(6, 10).Value = 0
$ws1.Cells.Item(6, 11).Value = 0.06376606877792257
$ws1.Cells.Item(6, 12).Value = 0
$ws1.Cells.Item(6, 13).Value = 0.01848318056051693
$ws1.Cells.Item(6, 14).Value = 0
$ws1.Cells.Item(6, 15).Value = 0
$ws1.Cells.Item(6, 16).Value = 0.02013165879070142
$ws1.Cells.Item(6, 17).Value = 0
$ws1.Cells.Item(6, 18).Value = 0
$ws1.Cells.Item(6, 19).Value = 0
$ws1.Cells.Item(6, 20).Value = 0
$ws1.Cells.Item(6, 21).Value = 0
$ws1.Cells.Item(6, 22).Value = 0
$ws1.Cells.Item(6, 23).Value = 0.02091671466191082
$ws1.Cells.Item(6, 24).Value = 0
$ws1.Cells.Item(6, 25).Value = 0.01092779432766532
$ws1.Cells.Item(6, 26).Value = 0
$ws1.Cells.Item(6, 27).Value = 0.02044461078705474
$ws1.Cells.Item(6, 28).Value = 0.01829069457430881
$ws1.Cells.Item(6, 29).Value = 0.0158530554669128
$ws1.Cells.Item(6, 30).Value = 0
$ws1.Cells.Item(6, 31).Value = 0.00958035760792172
$ws1.Cells.Item(6, 32).Value = 0
$ws1.Cells.Item(6, 33).Value = 0.01028201852480468
$ws1.Cells.Item(6, 34).Value = 0
$ws1.Cells.Item(6, 35).Value = 0.0

# --- Step2_Sj sheet: cumulative probability values ---
$ws2 = $wb.Worksheets.Item("Step2_Sj")
$ws2.Cells.Item(2, 4).Value = 0.05027127526523962
$ws2.Cells.Item(2, 5).Value = 0.1725294555728262
$ws2.Cells.Item(2, 6).Value = 0.6189759201403591
$ws2.Cells.Item(2, 7).Value = 0.7410848598483128
$ws2.Cells.Item(2, 8).Value = 0.7580863349604828
$ws2.Cells.Item(2, 9).Value = 0.7724075980785629
$ws2.Cells.Item(2, 10).Value = 0.7724075980785629
$ws2.Cells.Item(2, 11).Value = 0.7789372582241474
$ws2.Cells.Item(2, 12).Value = 0.7800197523633031
$ws2.Cells.Item(2, 13).Value = 0.7800197523633031
$ws2.Cells.Item(2, 14).Value = 0.7809327154485978
$ws2.Cells.Item(2, 15).Value = 0.7883247998725219
$ws2.Cells.Item(2, 16).Value = 0.8428309142500058
$ws2.Cells.Item(2, 17).Value = 0.8428309142500058
$ws2.Cells.Item(2, 18).Value = 0.8428309142500058
$ws2.Cells.Item(2, 19).Value = 0.8428309142500058
$ws2.Cells.Item(2, 20).Value = 0.8428309142500058
$ws2.Cells.Item(2, 21).Value = 0.8428309142500058
$ws2.Cells.Item(2, 22).Value = 0.8428309142500058
$ws2.Cells.Item(2, 23).Value = 0.8688458328274479
$ws2.Cells.Item(2, 24).Value = 0.8688458328274479
$ws2.Cells.Item(2, 25).Value = 0.8889208971979826
$ws2.Cells.Item(2, 26).Value = 0.8889208971979826
$ws2.Cells.Item(2, 27).Value = 0.9276221480863359
$ws2.Cells.Item(2, 28).Value = 0.9411253673983557
$ws2.Cells.Item(2, 29).Value = 0.9614418215839133
$ws2.Cells.Item(2, 30).Value = 0.9614418215839133
$ws2.Cells.Item(2, 31).Value = 0.9883505418514178
$ws2.Cells.Item(2, 32).Value = 0.9883505418514178
$ws2.Cells.Item(2, 33).Value = 1
$ws2.Cells.Item(2, 34).Value = 1
$ws2.Cells.Item(2, 35).Value = 1
$ws2.Cells.Item(3, 4).Value = 0.05026589338730262
$ws2.Cells.Item(3, 5).Value = 0.2093206633054425
$ws2.Cells.Item(3, 6).Value = 0.5656386447306766
$ws2.Cells.Item(3, 7).Value = 0.7401908915175625
$ws2.Cells.Item(3, 8).Value = 0.7626533630554873
$ws2.Cells.Item(3, 9).Value = 0.7846564305319357
$ws2.Cells.Item(3, 10).Value = 0.7846564305319357
$ws2.Cells.Item(3, 11).Value = 0.7926399432703481
$ws2.Cells.Item(3, 12).Value = 0.7926399432703481
$ws2.Cells.Item(3, 13).Value = 0.7926399432703481
$ws2.Cells.Item(3, 14).Value = 0.7926399432703481
$ws2.Cells.Item(3, 15).Value = 0.7926399432703481
$ws2.Cells.Item(3, 16).Value = 0.8206384082757682
$ws2.Cells.Item(3, 17).Value = 0.8206384082757682
$ws2.Cells.Item(3, 18).Value = 0.8206384082757682
$ws2.Cells.Item(3, 19).Value = 0.8206384082757682
$ws2.Cells.Item(3, 20).Value = 0.8206384082757682
$ws2.Cells.Item(3, 21).Value = 0.8206384082757682
$ws2.Cells.Item(3, 22).Value = 0.8206384082757682
$ws2.Cells.Item(3, 23).Value = 0.8503818499862704
$ws2.Cells.Item(3, 24).Value = 0.8503818499862704
$ws2.Cells.Item(3, 25).Value = 0.8754495840702333
$ws2.Cells.Item(3, 26).Value = 0.8754495840702333
$ws2.Cells.Item(3, 27).Value = 0.9120871831958848
$ws2.Cells.Item(3, 28).Value = 0.9264428190876398
$ws2.Cells.Item(3, 29).Value = 0.9479440760229466
$ws2.Cells.Item(3, 30).Value = 0.9479440760229466
$ws2.Cells.Item(3, 31).Value = 0.9828202508171013
$ws2.Cells.Item(3, 32).Value = 0.9828202508171013
$ws2.Cells.Item(3, 33).Value = 1
$ws2.Cells.Item(3, 34).Value = 1
$ws2.Cells.Item(3, 35).Value = 1
$ws2.Cells.Item(4, 4).Value = 0.2943213283311363
$ws2.Cells.Item(4, 5).Value = 0.3466231449229466
$ws2.Cells.Item(4, 6).Value = 0.7063919338685862
$ws2.Cells.Item(4, 7).Value = 0.7063919338685862
$ws2.Cells.Item(4, 8).Value = 0.7135755715906837
$ws2.Cells.Item(4, 9).Value = 0.7135755715906837
$ws2.Cells.Item(4, 10).Value = 0.7135755715906837
$ws2.Cells.Item(4, 11).Value = 0.7135755715906837
$ws2.Cells.Item(4, 12).Value = 0.7428372979286249
$ws2.Cells.Item(4, 13).Value = 0.7435892614247279
$ws2.Cells.Item(4, 14).Value = 0.7475525289544237
$ws2.Cells.Item(4, 15).Value = 0.8404173225029611
$ws2.Cells.Item(4, 16).Value = 0.868490453641917
$ws2.Cells.Item(4, 17).Value = 0.868490453641917
$ws2.Cells.Item(4, 18).Value = 0.868490453641917
$ws2.Cells.Item(4, 19).Value = 0.868490453641917
$ws2.Cells.Item(4, 20).Value = 0.868490453641917
$ws2.Cells.Item(4, 21).Value = 0.868490453641917
$ws2.Cells.Item(4, 22).Value = 0.9248133525077032
$ws2.Cells.Item(4, 23).Value = 0.9248133525077032
$ws2.Cells.Item(4, 24).Value = 0.9248133525077032
$ws2.Cells.Item(4, 25).Value = 0.9248133525077032
$ws2.Cells.Item(4, 26).Value = 0.9248133525077032
$ws2.Cells.Item(4, 27).Value = 0.9637954890644226
$ws2.Cells.Item(4, 28).Value = 0.9637954890644226
$ws2.Cells.Item(4, 29).Value = 0.9691014605956813
$ws2.Cells.Item(4, 30).Value = 0.9821435538032542
$ws2.Cells.Item(4, 31).Value = 0.9821435538032542
$ws2.Cells.Item(4, 32).Value = 0.9958895198621707
$ws2.Cells.Item(4, 33).Value = 1
$ws2.Cells.Item(4, 34).Value = 1
$ws2.Cells.Item(4, 35).Value = 1
$ws2.Cells.Item(5, 4).Value = 0.1729002386830045
$ws2.Cells.Item(5, 5).Value = 0.180549022868768
$ws2.Cells.Item(5, 6).Value = 0.7375528362879883
$ws2.Cells.Item(5, 7).Value = 0.78029758611869
$ws2.Cells.Item(5, 8).Value = 0.8216713910193975
$ws2.Cells.Item(5, 9).Value = 0.8258163196767285
$ws2.Cells.Item(5, 10).Value = 0.8258163196767285
$ws2.Cells.Item(5, 11).Value = 0.8443608081741073
$ws2.Cells.Item(5, 12).Value = 0.8443608081741073
$ws2.Cells.Item(5, 13).Value = 0.8443608081741073
$ws2.Cells.Item(5, 14).Value = 0.8443608081741073
$ws2.Cells.Item(5, 15).Value = 0.8685607769309016
$ws2.Cells.Item(5, 16).Value = 0.8983533950447625
$ws2.Cells.Item(5, 17).Value = 0.8983533950447625
$ws2.Cells.Item(5, 18).Value = 0.8983533950447625
$ws2.Cells.Item(5, 19).Value = 0.8983533950447625
$ws2.Cells.Item(5, 20).Value = 0.8983533950447625
$ws2.Cells.Item(5, 21).Value = 0.8983533950447625
$ws2.Cells.Item(5, 22).Value = 0.8998277901271146
$ws2.Cells.Item(5, 23).Value = 0.9158321600470107
$ws2.Cells.Item(5, 24).Value = 0.9158321600470107
$ws2.Cells.Item(5, 25).Value = 0.9312755481092244
$ws2.Cells.Item(5, 26).Value = 0.9312755481092244
$ws2.Cells.Item(5, 27).Value = 0.9674275861609367
$ws2.Cells.Item(5, 28).Value = 0.9714218741428857
$ws2.Cells.Item(5, 29).Value = 0.9900997006261018
$ws2.Cells.Item(5, 30).Value = 0.9900997006261018
$ws2.Cells.Item(5, 31).Value = 1
$ws2.Cells.Item(5, 32).Value = 1.0
$ws2.Cells.Item(5, 33).Value = 1.0
$ws2.Cells.Item(5, 34).Value = 1.0
$ws2.Cells.Item(5, 35).Value = 1.0
$ws2.Cells.Item(6, 4).Value = 0.03304963354092785
$ws2.Cells.Item(6, 5).Value = 0.2044518627093985
$ws2.Cells.Item(6, 6).Value = 0.5457098796379127
$ws2.Cells.Item(6, 7).Value = 0.711493662029204
$ws2.Cells.Item(6, 8).Value = 0.7582576829327287
$ws2.Cells.Item(6, 9).Value = 0.7913238459202805
$ws2.Cells.Item(6, 10).Value = 0.7913238459202805
$ws2.Cells.Item(6, 11).Value = 0.8550899146982031
$ws2.Cells.Item(6, 12).Value = 0.8550899146982031
$ws2.Cells.Item(6, 13).Value = 0.87357309525872
$ws2.Cells.Item(6, 14).Value = 0.87357309525872
$ws2.Cells.Item(6, 15).Value = 0.87357309525872
$ws2.Cells.Item(6, 16).Value = 0.8937047540494214
$ws2.Cells.Item(6, 17).Value = 0.8937047540494214
$ws2.Cells.Item(6, 18).Value = 0.8937047540494214
$ws2.Cells.Item(6, 19).Value = 0.8937047540494214
$ws2.Cells.Item(6, 20).Value = 0.8937047540494214
$ws2.Cells.Item(6, 21).Value = 0.8937047540494214
$ws2.Cells.Item(6, 22).Value = 0.8937047540494214
$ws2.Cells.Item(6, 23).Value = 0.9146214687113322
$ws2.Cells.Item(6, 24).Value = 0.9146214687113322
$ws2.Cells.Item(6, 25).Value = 0.9255492630389975
$ws2.Cells.Item(6, 26).Value = 0.9255492630389975
$ws2.Cells.Item(6, 27).Value = 0.9459938738260523
$ws2.Cells.Item(6, 28).Value = 0.964284568400361
$ws2.Cells.Item(6, 29).Value = 0.9801376238672739
$ws2.Cells.Item(6, 30).Value = 0.9801376238672739
$ws2.Cells.Item(6, 31).Value = 0.9897179814751956
$ws2.Cells.Item(6, 32).Value = 0.9897179814751956
$ws2.Cells.Item(6, 33).Value = 1
$ws2.Cells.Item(6, 34).Value = 1.0
$ws2.Cells.Item(6, 35).Value = 1.0
# --- Step3_DataPts_* sheets: tire-type filtering stat changes ---
$ws3 = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws3.Cells.Item(2, 3).Value = 1
$ws3.Cells.Item(2, 6).Value = 0.6189759201403591
$ws3.Cells.Item(2, 7).Value = 4
$ws3.Cells.Item(3, 3).Value = 1
$ws3.Cells.Item(3, 6).Value = 0.5656386447306766
$ws3.Cells.Item(3, 7).Value = 4
$ws3.Cells.Item(4, 6).Value = 0.7063919338685862
$ws3.Cells.Item(5, 6).Value = 0.7375528362879883
$ws3.Cells.Item(6, 6).Value = 0.5457098796379127

$ws4 = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws4.Cells.Item(2, 3).Value = 1
$ws4.Cells.Item(2, 4).Value = 6
$ws4.Cells.Item(2, 6).Value = 0.7410848598483128
$ws4.Cells.Item(2, 7).Value = 5
$ws4.Cells.Item(3, 3).Value = 1
$ws4.Cells.Item(3, 4).Value = 6
$ws4.Cells.Item(3, 6).Value = 0.7401908915175625
$ws4.Cells.Item(3, 7).Value = 5
$ws4.Cells.Item(4, 6).Value = 0.7063919338685862
$ws4.Cells.Item(5, 6).Value = 0.7375528362879883
$ws4.Cells.Item(6, 6).Value = 0.711493662029204

$ws5 = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws5.Cells.Item(2, 3).Value = 1
$ws5.Cells.Item(2, 4).Value = 15
$ws5.Cells.Item(2, 6).Value = 0.8428309142500058
$ws5.Cells.Item(2, 7).Value = 14
$ws5.Cells.Item(3, 3).Value = 1
$ws5.Cells.Item(3, 4).Value = 15
$ws5.Cells.Item(3, 6).Value = 0.8206384082757682
$ws5.Cells.Item(3, 7).Value = 14
$ws5.Cells.Item(4, 4).Value = 14
$ws5.Cells.Item(4, 6).Value = 0.8404173225029611
$ws5.Cells.Item(4, 7).Value = 13
$ws5.Cells.Item(5, 4).Value = 7
$ws5.Cells.Item(5, 6).Value = 0.8216713910193975
$ws5.Cells.Item(5, 7).Value = 6
$ws5.Cells.Item(6, 4).Value = 10
$ws5.Cells.Item(6, 6).Value = 0.8550899146982031
$ws5.Cells.Item(6, 7).Value = 8

$ws6 = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws6.Cells.Item(2, 3).Value = 1
$ws6.Cells.Item(2, 4).Value = 26
$ws6.Cells.Item(2, 6).Value = 0.9276221480863359
$ws6.Cells.Item(2, 7).Value = 25
$ws6.Cells.Item(3, 3).Value = 1
$ws6.Cells.Item(3, 4).Value = 26
$ws6.Cells.Item(3, 6).Value = 0.9120871831958848
$ws6.Cells.Item(3, 7).Value = 25
$ws6.Cells.Item(4, 4).Value = 21
$ws6.Cells.Item(4, 6).Value = 0.9248133525077032
$ws6.Cells.Item(4, 7).Value = 20
$ws6.Cells.Item(5, 4).Value = 22
$ws6.Cells.Item(5, 6).Value = 0.9158321600470107
$ws6.Cells.Item(5, 7).Value = 21
$ws6.Cells.Item(6, 4).Value = 22
$ws6.Cells.Item(6, 6).Value = 0.9146214687113322
$ws6.Cells.Item(6, 7).Value = 20

# --- Tire_Type shared text fix across all Step3_DataPts_* sheets (K2:K6) ---
$ws3.Range("K2:K6").Value = "710R"
$ws4.Range("K2:K6").Value = "710R"
$ws5.Range("K2:K6").Value = "710R"
$ws6.Range("K2:K6").Value = "710R"